$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old A1:B5 block (the date/amount pair in row 1 plus the
# trailing empty rows 2-5) - it is being relocated down to row 7.
$ws.Range("A1:B5").Clear()

# A7 holds the date as literal text ("2023-07-18"), not an auto-converted
# serial date. Enter it as a formula returning that literal text, then
# paste-special as values-only so the stored cell is plain text (no
# formula, no extra number-format/quote-prefix styling) - exactly like a
# user typing the formula in and then "Paste Values" over it.
$ws.Range("A7").Formula = "=""2023-07-18"""
$ws.Range("A7").Copy()
$ws.Range("A7").PasteSpecial(-4163)
$excel.CutCopyMode = $false

# B7 holds the accompanying numeric amount.
$ws.Range("B7").Value = 25442
